# TanLoc_Commit11: Update dong goi api
# Renumber the EPC (col B) and Serial (col D) values for the 130 data rows
# (rows 2-131) so the sequence restarts at 30300E890A0180C077359401 /
# 2000000001 instead of 30300E890A0180C077359483 / 2000000131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base values corresponding to data row 2 (the first product row).
$baseEpc    = [System.Convert]::ToInt64("30300E890A0180C077359401".Substring(16), 16)
$epcPrefix  = "30300E890A0180C0"
$baseSerial = 2000000001

$firstRow = 2
$lastRow  = 131

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $offset = $row - $firstRow

    $epcSuffix = "{0:X8}" -f ($baseEpc + $offset)
    $epcValue  = $epcPrefix + $epcSuffix

    $serialValue = $baseSerial + $offset

    # EPC is hexadecimal text (never parses as a plain number) so a plain
    # string assignment keeps it stored as text.
    $ws.Cells.Item($row, 2).Value = $epcValue

    # Serial looks like a pure integer, so prefix with an apostrophe to
    # force Excel to keep storing it as text (matching the source file).
    $ws.Cells.Item($row, 4).Value = "'" + $serialValue
}
